# "Still refactoring the prompt classes"
# - Rename the "sample_prompts" sheet to "prompts"
# - Update cell C25 on that sheet from 6 to 2
# - Move the sheet's active selection from H28 to C25 (and drop the
#   scrolled-away topLeftCell view state, since the selected cell is
#   now back in the default top-left viewport)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample_prompts")

# Rename sheet
$ws.Name = "prompts"

# Update the data value
$ws.Range("C25").Value = 2

# Update the view: make this sheet active and select C25
$ws.Activate() | Out-Null
$ws.Range("C25").Select() | Out-Null
